$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look numeric but must remain plain text (match original inlineStr formatting).
# We temporarily force a Text number format before assignment, then restore the default "Normal" style
# so the cell keeps no explicit style (matching the original workbook which has no style on these cells).
$textCells = @{
    'D5' = '245.82'
    'D6' = '0.629'
    'D7' = '74.62'
    'D9' = '0.627'
    'D10' = '39.68'
    'D11' = '0.0944'
    'D12' = '7.14'
    'D13' = '0.104'
    'D15' = '14.77'
    'D16' = '0.857'
    'D20' = '6.11'
    'D21' = '71.29'
    'D23' = '230.30'
    'D24' = '0.999'
    'D25' = '3.71'
    'D26' = '11.23'
    'D27' = '2.30'
    'D28' = '7.21'
    'D30' = '169.21'
    'D31' = '20.50'
    'D32' = '34.31'
    'D33' = '0.0838'
    'D34' = '0.119'
    'D36' = '4.63'
    'D37' = '4.88'
    'D38' = '0.0298'
    'D39' = '13.52'
    'D40' = '5.89'
    'D41' = '2.18'
    'D42' = '110.33'
    'D43' = '0.202'
    'D44' = '60.12'
    'D45' = '8.76'
    'D47' = '0.996'
    'D50' = '4.24'
    'D51' = '2.24'
}
foreach ($addr in $textCells.Keys) {
    $ws.Range($addr).NumberFormat = "@"
}
foreach ($addr in $textCells.Keys) {
    $ws.Range($addr).Value = $textCells[$addr]
}
foreach ($addr in $textCells.Keys) {
    $ws.Range($addr).Style = "Normal"
}

# Remaining cells are safe to set directly (non-numeric-looking text, percent strings, etc.)
$ws.Range("D2").Value = "41.838.47"
$ws.Range("E2").Value = "  -1.89%  "
$ws.Range("D3").Value = "2.233.33"
$ws.Range("E3").Value = "  -2.18%  "
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("E5").Value = "  -2.40%  "
$ws.Range("E6").Value = "  -1.20%  "
$ws.Range("E7").Value = "  +1.19%  "
$ws.Range("E8").Value = "  +0.12%  "
$ws.Range("E9").Value = "  -3.15%  "
$ws.Range("E10").Value = "  +1.36%  "
$ws.Range("E11").Value = "  -4.75%  "
$ws.Range("E12").Value = "  -2.89%  "
$ws.Range("E13").Value = "  -1.81%  "
$ws.Range("D14").Value = "2.578.75"
$ws.Range("E14").Value = "  -1.77%  "
$ws.Range("E15").Value = "  -1.38%  "
$ws.Range("E16").Value = "  -2.02%  "
$ws.Range("D17").Value = "2.225.19"
$ws.Range("E17").Value = "  -2.26%  "
$ws.Range("D18").Value = "41.774.13"
$ws.Range("E18").Value = "  -1.89%  "
$ws.Range("D19").Value = "0.0₃0976"
$ws.Range("E19").Value = "  -3.15%  "
$ws.Range("E20").Value = "  -3.24%  "
$ws.Range("E21").Value = "  -1.44%  "
$ws.Range("E22").Value = "  +1.63%  "
$ws.Range("E23").Value = "  -2.00%  "
$ws.Range("E24").Value = "  -0.03%  "
$ws.Range("B25").Value = "WEMIXToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("E25").Value = "  -5.15%  "
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("E26").Value = "  -2.18%  "
$ws.Range("E27").Value = "  -5.11%  "
$ws.Range("E28").Value = "  +13.17%  "
$ws.Range("E29").Value = "  -1.46%  "
$ws.Range("E30").Value = "  +1.39%  "
$ws.Range("E31").Value = "  -2.76%  "
$ws.Range("E32").Value = "  +6.43%  "
$ws.Range("E33").Value = "  +2.29%  "
$ws.Range("E34").Value = "  -6.62%  "
$ws.Range("E35").Value = "  +0.24%  "
$ws.Range("E36").Value = "  -2.11%  "
$ws.Range("E37").Value = "  +2.69%  "
$ws.Range("E38").Value = "  -2.73%  "
$ws.Range("E39").Value = "  -2.40%  "
$ws.Range("E40").Value = "  -1.24%  "
$ws.Range("E41").Value = "  -6.77%  "
$ws.Range("E42").Value = "  +13.20%  "
$ws.Range("E43").Value = "  -5.53%  "
$ws.Range("E44").Value = "  -2.98%  "
$ws.Range("E45").Value = "  -4.25%  "
$ws.Range("E46").Value = "  -3.25%  "
$ws.Range("E47").Value = "  -0.42%  "
$ws.Range("E48").Value = "  -4.46%  "
$ws.Range("E49").Value = "  -1.68%  "
$ws.Range("E50").Value = "  -12.14%  "
$ws.Range("E51").Value = "  -2.49%  "
